# Insert a new data row at row 121 (shifts existing rows 121:232 down to 122:233)
# and populate it with the new record's values, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 121..232 down by one row, creating a blank row 121.
$ws.Rows.Item(121).Insert()

# Populate the newly inserted row 121 with the new record.
$ws.Cells.Item(121, 1).Value  = 11
$ws.Cells.Item(121, 2).Value  = 'Vega Monumental Concepción'
$ws.Cells.Item(121, 3).Value  = 'Bíobío'
$ws.Cells.Item(121, 4).Value  = 44904
$ws.Cells.Item(121, 5).Value  = 8
$ws.Cells.Item(121, 6).Value  = 100112003
$ws.Cells.Item(121, 7).Value  = 'Ajo'
$ws.Cells.Item(121, 8).Value  = 'Chino'
$ws.Cells.Item(121, 9).Value  = 'Primera'
$ws.Cells.Item(121, 10).Value = 220
$ws.Cells.Item(121, 11).Value = 13000
$ws.Cells.Item(121, 12).Value = 14000
$ws.Cells.Item(121, 13).Value = 13455
$ws.Cells.Item(121, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(121, 15).Value = 'China'
$ws.Cells.Item(121, 16).Value = 1346
$ws.Cells.Item(121, 17).Value = 10
$ws.Cells.Item(121, 18).Value = 'Hortaliza'
